$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F3: 10~12/L -> 10^12/L
$ws.Range("F3").Value = "10^12/L"

# F6: fl -> fL
$ws.Range("F6").Value = "fL"

# C13: 7.700 -> 7.7 (keep as text, not numeric)
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "7.7"
$ws.Range("C13").Style = "Normal"

# C14: 0.700 -> 0.7 (keep as text, not numeric)
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "0.7"
$ws.Range("C14").Style = "Normal"

# C15: 14.000 -> 14.0 (keep as text, not numeric)
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "14.0"
$ws.Range("C15").Style = "Normal"

# Row 16: split combined code/label into A16 (code) and B16 (label)
$ws.Range("A16").Value = "RDW-CV"
$ws.Range("B16").Value = "红细胞平均大小"

# Row 17: split combined code/label into A17 (code) and B17 (label)
$ws.Range("A17").Value = "RDW-SD"
$ws.Range("B17").Value = "红细胞平均宽度"

# C18: 8.500 -> 8.5 (keep as text, not numeric)
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "8.5"
$ws.Range("C18").Style = "Normal"

# Row 21: clear A21 (becomes blank, cell kept), set B21 to RH血型
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = ""
$ws.Range("A21").Style = "Normal"
$ws.Range("B21").Value = "RH血型"
